$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.283.00'
$ws.Range("E2").Value = '  +0.45%  '

$ws.Range("D3").Value = '2.083.83'
$ws.Range("E3").Value = '  +0.23%  '

$ws.Range("E4").Value = '  -0.20%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '252.03'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.677'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +4.38%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '62.61'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +25.02%  '

$ws.Range("E8").Value = '  -0.16%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '61.87'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.43%  '

$ws.Range("E10").Value = '  +5.61%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0817'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +10.92%  '

$ws.Range("E12").Value = '  +2.97%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '15.71'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +3.89%  '

$ws.Range("D14").Value = '2.377.46'
$ws.Range("E14").Value = '  -0.26%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.830'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.50%  '

$ws.Range("E16").Value = '  +7.55%  '

$ws.Range("D17").Value = '2.082.80'
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("D18").Value = '37.184.37'
$ws.Range("E18").Value = '  +0.68%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '74.92'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +4.22%  '

$ws.Range("D20").Value = '0.0₃0931'
$ws.Range("E20").Value = '  +13.93%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '15.19'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +15.65%  '

$ws.Range("E22").Value = '  +5.92%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '240.53'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.16%  '

$ws.Range("E24").Value = '  +0.04%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.82%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '171.71'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.86%  '

$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '20.54'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.57%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.06'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +3.79%  '

$ws.Range("E30").Value = '  +3.75%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.82'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +8.08%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.09'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +2.60%  '

$ws.Range("E33").Value = '  +6.06%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.47'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +10.84%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0893'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.29%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.33'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +3.10%  '

$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("E38").Value = '  -3.48%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.110'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +24.74%  '

$ws.Range("E40").Value = '  +3.32%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '18.85'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +8.05%  '

$ws.Range("E42").Value = '  +3.03%  '

$ws.Range("E43").Value = '  +2.65%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '99.29'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.21%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '4.44'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +25.28%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.81'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.39%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +14.50%  '

$ws.Range("E48").Value = '  +15.76%  '

$ws.Range("D49").Value = '1.313.38'
$ws.Range("E49").Value = '  +0.95%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.94'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.18%  '
